$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) - reuse the same header style (s="1") already
# used by the other header cells (B1:G1) by copying formats from G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" data values for the two data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
